$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value2 = 2.76
$ws.Range("G2").Value2 = 2.86
$ws.Range("H2").Value2 = 2.5
$ws.Range("I2").Value2 = 2.58
$ws.Range("J2").Value2 = 3.95
$ws.Range("K2").Value2 = 4.1
$ws.Range("M2").Value2 = 1.03
$ws.Range("N2").Value2 = 5.8
$ws.Range("O2").Value2 = 1.19
$ws.Range("Q2").Value2 = 1.59
$ws.Range("R2").Value2 = 1.68
$ws.Range("S2").Value2 = 2.44
$ws.Range("U2").Value2 = 2.66
$ws.Range("V2").Value2 = 1.63
$ws.Range("W2").Value2 = 1.53
$ws.Range("Y2").Value2 = 19
$ws.Range("Z2").Value2 = 23
$ws.Range("AA2").Value2 = 65
$ws.Range("AC2").Value2 = 9.800000000000001
$ws.Range("AE2").Value2 = 23
$ws.Range("AF2").Value2 = 26
$ws.Range("AH2").Value2 = 15
$ws.Range("AI2").Value2 = 48
$ws.Range("AJ2").Value2 = 55
$ws.Range("AL2").Value2 = 48
$ws.Range("AM2").Value2 = 80
$ws.Range("AN2").Value2 = 19
$ws.Range("AO2").Value2 = 13.5
$ws.Range("F3").Value2 = 1.48
$ws.Range("G3").Value2 = 1.58
$ws.Range("H3").Value2 = 4.9
$ws.Range("I3").Value2 = 6.6
$ws.Range("J3").Value2 = 5.3
$ws.Range("K3").Value2 = 6
$ws.Range("L3").Value2 = 1.17
$ws.Range("N3").Value2 = 1.11
$ws.Range("O3").Value2 = 1.1
$ws.Range("P3").Value2 = 3.1
$ws.Range("Q3").Value2 = 1.33
$ws.Range("R3").Value2 = 2
$ws.Range("S3").Value2 = 1.8
$ws.Range("T3").Value2 = 1.48
$ws.Range("U3").Value2 = 2.7
$ws.Range("V3").Value2 = 1.18
$ws.Range("W3").Value2 = 2.68
$ws.Range("X3").Value2 = 1000
$ws.Range("Y3").Value2 = 90
$ws.Range("Z3").Value2 = 70
$ws.Range("AA3").Value2 = 830
$ws.Range("AB3").Value2 = 19.5
$ws.Range("AE3").Value2 = 65
$ws.Range("AF3").Value2 = 1000
$ws.Range("AG3").Value2 = 11.5
$ws.Range("AH3").Value2 = 23
$ws.Range("AI3").Value2 = 180
$ws.Range("AJ3").Value2 = 18.5
$ws.Range("AK3").Value2 = 14
$ws.Range("AL3").Value2 = 22
$ws.Range("AN3").Value2 = 6
$ws.Range("AO3").Value2 = 1000
$ws.Range("F4").Value2 = 1.51
$ws.Range("G4").Value2 = 1.57
$ws.Range("I4").Value2 = 8.6
$ws.Range("J4").Value2 = 4.1
$ws.Range("U4").Value2 = 1.82
$ws.Range("W4").Value2 = 2.74
$ws.Range("Y4").Value2 = 26
$ws.Range("AB4").Value2 = 8.6
$ws.Range("AC4").Value2 = 11.5
$ws.Range("AF4").Value2 = 10
$ws.Range("AG4").Value2 = 12
$ws.Range("AJ4").Value2 = 16.5
$ws.Range("AK4").Value2 = 20
$ws.Range("AL4").Value2 = 50
$ws.Range("AN4").Value2 = 10.5
$ws.Range("H5").Value2 = 2.38
$ws.Range("J5").Value2 = 2.84
$ws.Range("K5").Value2 = 3.5
$ws.Range("M5").Value2 = 1.05
$ws.Range("O5").Value2 = 1.05
$ws.Range("S5").Value2 = 2.14
$ws.Range("P6").Value2 = 1.94
$ws.Range("F7").Value2 = 1.33
$ws.Range("I7").Value2 = 12.5
$ws.Range("J7").Value2 = 5.4
$ws.Range("K7").Value2 = 6.2
$ws.Range("L7").Value2 = 1.28
$ws.Range("M7").Value2 = 1.04
$ws.Range("N7").Value2 = 4.7
$ws.Range("O7").Value2 = 1.21
$ws.Range("P7").Value2 = 2.32
$ws.Range("Q7").Value2 = 1.63
$ws.Range("R7").Value2 = 1.5
$ws.Range("S7").Value2 = 2.58
$ws.Range("T7").Value2 = 2.04
$ws.Range("U7").Value2 = 1.81
$ws.Range("V7").Value2 = 1.09
$ws.Range("W7").Value2 = 3.6
$ws.Range("X7").Value2 = 27
$ws.Range("Y7").Value2 = 1000
$ws.Range("AF7").Value2 = 10
$ws.Range("AG7").Value2 = 11
$ws.Range("AH7").Value2 = 36
$ws.Range("AJ7").Value2 = 13
$ws.Range("AN7").Value2 = 5.8
$ws.Range("G8").Value2 = 2.92
$ws.Range("H8").Value2 = 2.58
$ws.Range("J8").Value2 = 3.6
$ws.Range("K8").Value2 = 3.7
$ws.Range("Q8").Value2 = 1.83
$ws.Range("S8").Value2 = 3.05
$ws.Range("T8").Value2 = 1.68
$ws.Range("U8").Value2 = 2.26
$ws.Range("V8").Value2 = 1.58
$ws.Range("W8").Value2 = 1.52
$ws.Range("X8").Value2 = 15.5
$ws.Range("Z8").Value2 = 18
$ws.Range("AE8").Value2 = 28
$ws.Range("AJ8").Value2 = 48
$ws.Range("AL8").Value2 = 46
$ws.Range("F9").Value2 = 4.1
$ws.Range("K9").Value2 = 4.4
$ws.Range("F10").Value2 = 2.36
$ws.Range("G10").Value2 = 2.48
$ws.Range("H10").Value2 = 3.35
$ws.Range("I10").Value2 = 3.65
$ws.Range("J10").Value2 = 3.3
$ws.Range("N10").Value2 = 3.45
$ws.Range("P10").Value2 = 1.84
$ws.Range("T10").Value2 = 1.77
$ws.Range("U10").Value2 = 2.08
$ws.Range("V10").Value2 = 1.38
$ws.Range("W10").Value2 = 1.67
$ws.Range("X10").Value2 = 13.5
$ws.Range("Y10").Value2 = 13
$ws.Range("Z10").Value2 = 24
$ws.Range("AA10").Value2 = 65
$ws.Range("AE10").Value2 = 42
$ws.Range("AF10").Value2 = 18
$ws.Range("AK10").Value2 = 32
$ws.Range("M11").Value2 = 1.02
$ws.Range("Q11").Value2 = 1.02
$ws.Range("G12").Value2 = 3.15
$ws.Range("H12").Value2 = 2.26
$ws.Range("I12").Value2 = 2.36
$ws.Range("K12").Value2 = 4.2
$ws.Range("R12").Value2 = 1.61
$ws.Range("T12").Value2 = 1.54
$ws.Range("V12").Value2 = 1.73
$ws.Range("W12").Value2 = 1.46
$ws.Range("AG12").Value2 = 1000
$ws.Range("F13").Value2 = 2.8
$ws.Range("I13").Value2 = 2.86
$ws.Range("O13").Value2 = 1.35
$ws.Range("Q13").Value2 = 2.02
$ws.Range("S13").Value2 = 3.6
$ws.Range("T13").Value2 = 1.79
$ws.Range("V13").Value2 = 1.54
$ws.Range("L14").Value2 = 1.43
$ws.Range("Q14").Value2 = 2.08
$ws.Range("U14").Value2 = 1.67
$ws.Range("F16").Value2 = 3.8
$ws.Range("G16").Value2 = 4.1
$ws.Range("H16").Value2 = 2.48
$ws.Range("I16").Value2 = 2.54
$ws.Range("J16").Value2 = 2.88
$ws.Range("K16").Value2 = 2.98
$ws.Range("P16").Value2 = 1.36
$ws.Range("V16").Value2 = 1.64
$ws.Range("W16").Value2 = 1.33
$ws.Range("X16").Value2 = 6.4
$ws.Range("Z16").Value2 = 13
$ws.Range("AA16").Value2 = 40
$ws.Range("AB16").Value2 = 9
$ws.Range("AC16").Value2 = 7.6
$ws.Range("AD16").Value2 = 14.5
$ws.Range("AE16").Value2 = 48
$ws.Range("AF16").Value2 = 25
$ws.Range("AH16").Value2 = 34
$ws.Range("AI16").Value2 = 110
$ws.Range("AJ16").Value2 = 100
$ws.Range("AM16").Value2 = 410
$ws.Range("AN16").Value2 = 160
$ws.Range("AO16").Value2 = 60
$ws.Range("F17").Value2 = 2
$ws.Range("G17").Value2 = 2.16
$ws.Range("I17").Value2 = 5.2
$ws.Range("N17").Value2 = 2.54
$ws.Range("O17").Value2 = 1.54
$ws.Range("T17").Value2 = 2.18
$ws.Range("U17").Value2 = 1.71
$ws.Range("AH17").Value2 = 30
$ws.Range("F18").Value2 = 4.3
$ws.Range("G18").Value2 = 5.7
$ws.Range("H18").Value2 = 1.91
$ws.Range("I18").Value2 = 2.18
$ws.Range("K18").Value2 = 3.95
$ws.Range("O18").Value2 = 1.54
$ws.Range("P18").Value2 = 1.5
$ws.Range("Q18").Value2 = 2.4
$ws.Range("S18").Value2 = 5.5
$ws.Range("T18").Value2 = 2.2
$ws.Range("U18").Value2 = 1.65
$ws.Range("V18").Value2 = 1.84
$ws.Range("W18").Value2 = 1.22
$ws.Range("X18").Value2 = 10.5
$ws.Range("Z18").Value2 = 12
$ws.Range("AB18").Value2 = 14.5
$ws.Range("AC18").Value2 = 9.6
$ws.Range("AG18").Value2 = 27
$ws.Range("I19").Value2 = 5.7
$ws.Range("L19").Value2 = 1.46
$ws.Range("N19").Value2 = 3.4
$ws.Range("P19").Value2 = 1.82
$ws.Range("Q19").Value2 = 2.18
$ws.Range("R19").Value2 = 1.32
$ws.Range("S19").Value2 = 4
$ws.Range("F20").Value2 = 3.6
$ws.Range("I20").Value2 = 2.36
$ws.Range("K20").Value2 = 3.5
$ws.Range("N20").Value2 = 2.96
$ws.Range("Q20").Value2 = 2.28
$ws.Range("S20").Value2 = 3.85
$ws.Range("T20").Value2 = 1.96
$ws.Range("V20").Value2 = 1.73
